# Applies the diff: swap the data of rows 4 and 5, and update B16.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# --- Swap simple valued cells between row 4 and row 5 ---
$cols = @("A","B","E","F","G","H","Q","R","AC","AE")
foreach ($col in $cols) {
    $addr4 = "$col" + "4"
    $addr5 = "$col" + "5"
    $v4 = $ws.Range($addr4).Value2
    $v5 = $ws.Range($addr5).Value2
    $ws.Range($addr4).Value2 = $v5
    $ws.Range($addr5).Value2 = $v4
}

# --- Handle J: row4 had "fruktkroppar", row5 was blank/absent ---
$jVal = $ws.Range("J4").Value2
$ws.Range("J5").Value2 = $jVal
$ws.Range("J4").ClearContents()

# --- Handle AF: row4 had an empty placeholder cell, row5 was absent ---
# Make AF5 a present-but-empty cell (copy from an existing empty placeholder cell),
# then clear AF4 so it becomes absent.
$ws.Range("K5").Copy($ws.Range("AF5"))
$ws.Range("AF4").ClearContents()

# --- Handle L and M: row5 had empty placeholder cells, row4 was absent ---
# Make L4 and M4 present-but-empty cells, then clear L5 and M5 so they become absent.
$ws.Range("K5").Copy($ws.Range("L4"))
$ws.Range("K5").Copy($ws.Range("M4"))
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()

# --- Update B16 ---
$ws.Range("B16").Value2 = 91834
